# Adds the "ODI Bowling Extra" worksheet (extra bowling scraping attributes)
# as the 5th sheet in the workbook, mirroring the structure of the existing
# "ODI Batting Extra" sheet.

$wb = $excel.ActiveWorkbook

# Use the existing "ODI Batting Extra" sheet as the anchor so the new sheet
# is inserted immediately after it (i.e. becomes the last / 5th sheet).
$anchor = $wb.Worksheets.Item("ODI Batting Extra")

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $anchor)
$newSheet.Name = "ODI Bowling Extra"

# Header row (same bold/centered/bordered header style used by the other sheets)
$newSheet.Cells.Item(1, 1).Value = "MATCH_CODE"
$newSheet.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$newSheet.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"

# Data rows (MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL).
# All values are written as text (matching the source scrape format, which
# stores even numeric-looking values such as match codes and percentages as
# plain strings). A leading apostrophe forces Excel to keep the literal text
# instead of auto-converting it to a number.
$data = @(
    ,@("4294", "", "")
    ,@("4300", "1", "10.00%")
    ,@("4426", "", "")
    ,@("4427", "1", "20.00%")
    ,@("4428", "1", "10.00%")
    ,@("4469", "1", "30.00%")
    ,@("4470", "1", "40.00%")
    ,@("4471", "", "")
    ,@("4598", "", "")
    ,@("4599", "", "")
    ,@("4602", "0", "40.00%")
    ,@("4609", "0", "")
    ,@("4613", "2", "10.00%")
    ,@("4618", "0", "")
    ,@("4620", "1", "10.00%")
    ,@("4622", "", "")
    ,@("4660", "", "")
    ,@("4663", "", "")
    ,@("4666", "", "")
    ,@("4698", "0", "")
)

$row = 2
foreach ($entry in $data) {
    $newSheet.Cells.Item($row, 1).Value = "'" + $entry[0]
    if ($entry[1] -ne "") {
        $newSheet.Cells.Item($row, 2).Value = "'" + $entry[1]
    }
    if ($entry[2] -ne "") {
        $newSheet.Cells.Item($row, 3).Value = "'" + $entry[2]
    }
    $row = $row + 1
}

$wb.Save()
